$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in shared string used by B23 ("unlick" -> "unlock").
$ws.Range("B23").Value = "unlock home doors"

# Fill in the previously-empty B column for rows 28-32.
$ws.Range("B28").Value = "home lights on"
$ws.Range("B29").Value = "thermostat home"
$ws.Range("B30").Value = "close garage door"
$ws.Range("B31").Value = "open garage door"
$ws.Range("B32").Value = "open garage door"

# Update the sheet view: drop the frozen/scrolled topLeftCell and move the
# active selection from H35 to E13.
$ws.Range("E13").Select()
